$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21. This pushes the existing rows 21..151
# down to 22..152 (and therefore the sheet dimension grows from R151 to R152).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new data record.
$ws.Cells.Item(21, 1).Value = 4
$ws.Cells.Item(21, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(21, 3).Value = "Los Lagos"
$ws.Cells.Item(21, 4).Value = 44473
$ws.Cells.Item(21, 5).Value = 10
$ws.Cells.Item(21, 6).Value = 100112043
$ws.Cells.Item(21, 7).Value = "Pepino ensalada"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 150
$ws.Cells.Item(21, 11).Value = 23000
$ws.Cells.Item(21, 12).Value = 23000
$ws.Cells.Item(21, 13).Value = 23000
$ws.Cells.Item(21, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(21, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value = 383
$ws.Cells.Item(21, 17).Value = 60
$ws.Cells.Item(21, 18).Value = "Hortaliza"
